# Update training metrics (rows 2-26) for the new LM training run.
# Every data row (2 through 26) shares the same metric values per column,
# so the new values are assigned column-by-column across the whole block.
# Note: numeric literals are written in plain decimal form (no scientific
# notation) because the script engine's lexer does not accept exponent
# suffixes like "e-06" on number literals.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    "B" = 0.9999969141016266
    "C" = 0.9990169585194482
    "D" = 0.9999895876125278
    "E" = 0.9999999199992772
    "F" = 0.9999972022415186
    "G" = 0.00000288055067624043
    "H" = 0.0009176260715613002
    "I" = 0.00000311858446964245
    "J" = 0.00000005281892399633883
    "K" = 0.000001585701696819395
    "L" = 0.00009999963621106398
    "M" = 0.00169721851163615
    "N" = 0.9999753128130129
    "O" = 0.001769472617403162
    "P" = 67.51505815013246
    "Q" = 93.11145047236465
}

foreach ($col in $newValues.Keys) {
    $range = $ws.Range("$($col)2:$($col)26")
    $range.Value = $newValues[$col]
}
